$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 927.75
$ws.Range("I43").Value = 870.5
$ws.Range("J43").Value = 985
$ws.Range("K43").Value = 870.5
$ws.Range("L43").Value = 985
$ws.Range("M43").Value = -801.5
$ws.Range("N43").Value = -1123
$ws.Range("H98").Value = 861705.0600000001
$ws.Range("I98").Value = 1118076.6
$ws.Range("J98").Value = 7133.3335
$ws.Range("K98").Value = 1118076.6
$ws.Range("L98").Value = 7133.3335
$ws.Range("M98").Value = -1116578.6
$ws.Range("N98").Value = -10129.3335
$ws.Range("H100").Value = 12823267
$ws.Range("I100").Value = 27779646
$ws.Range("J100").Value = 3514.2856
$ws.Range("K100").Value = 27779646
$ws.Range("L100").Value = 3514.2856
$ws.Range("M100").Value = -27779105
$ws.Range("N100").Value = -4596.2856
$ws.Range("H112").Value = 5815037
$ws.Range("J112").Value = 6411407.5
$ws.Range("L112").Value = 19234222.5
$ws.Range("N112").Value = -19236438.5
$ws.Range("H122").Value = 861705.0600000001
$ws.Range("I122").Value = 1118076.6
$ws.Range("J122").Value = 7133.3335
$ws.Range("K122").Value = 3354229.8
$ws.Range("L122").Value = 21400.0005
$ws.Range("M122").Value = -3351779.8
$ws.Range("N122").Value = -26300.0005
$ws.Range("H132").Value = 263929.1
$ws.Range("I132").Value = 419634.22
$ws.Range("J132").Value = 26273.947
$ws.Range("K132").Value = 1258902.66
$ws.Range("L132").Value = 78821.841
$ws.Range("M132").Value = -1256372.66
$ws.Range("N132").Value = -83881.841
$ws.Range("H138").Value = 6191267
$ws.Range("J138").Value = 8477519
$ws.Range("L138").Value = 25432557
$ws.Range("N138").Value = -25442837

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1897.6471
$ws.Range("I45").Value = 1828.75
$ws.Range("K45").Value = 1828.75
$ws.Range("M45").Value = -1451.75
$ws.Range("H122").Value = 1674.5834
$ws.Range("I122").Value = 1377
$ws.Range("K122").Value = 4131
$ws.Range("M122").Value = -1681
$ws.Range("H123").Value = 33618.332
$ws.Range("J123").Value = 33618.332
$ws.Range("L123").Value = 33618.332
$ws.Range("N123").Value = -43418.332
$ws.Range("H133").Value = 58959.2
$ws.Range("J133").Value = 58959.2
$ws.Range("L133").Value = 58959.2
$ws.Range("N133").Value = -64019.2
$ws.Range("H139").Value = 60000
$ws.Range("J139").Value = 60000
$ws.Range("L139").Value = 60000
$ws.Range("N139").Value = -70280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 126132
$ws.Range("I97").Value = 126132
$ws.Range("K97").Value = 126132
$ws.Range("M97").Value = -125141

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3565
$ws.Range("I31").Value = 1654.381
$ws.Range("J31").Value = 5309.478
$ws.Range("K31").Value = 1654.381
$ws.Range("L31").Value = 5309.478
$ws.Range("M31").Value = -1359.381
$ws.Range("N31").Value = -5899.478
$ws.Range("H34").Value = 3565
$ws.Range("I34").Value = 1654.381
$ws.Range("J34").Value = 5309.478
$ws.Range("K34").Value = 1654.381
$ws.Range("L34").Value = 5309.478
$ws.Range("M34").Value = -1452.381
$ws.Range("N34").Value = -5713.478
$ws.Range("H62").Value = 17325.934
$ws.Range("I62").Value = 27499.375
$ws.Range("J62").Value = 5699.143
$ws.Range("K62").Value = 27499.375
$ws.Range("L62").Value = 5699.143
$ws.Range("M62").Value = -26875.375
$ws.Range("N62").Value = -6947.143
$ws.Range("H65").Value = 17325.934
$ws.Range("I65").Value = 27499.375
$ws.Range("J65").Value = 5699.143
$ws.Range("K65").Value = 137496.875
$ws.Range("L65").Value = 28495.715
$ws.Range("M65").Value = -134376.875
$ws.Range("N65").Value = -34735.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 7500
$ws.Range("J104").Value = 7500
$ws.Range("L104").Value = 22500
$ws.Range("N104").Value = -27742
$ws.Range("H113").Value = 25001082
$ws.Range("I113").Value = 400
$ws.Range("J113").Value = 26316908
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 78950724
$ws.Range("M113").Value = 970
$ws.Range("N113").Value = -78955064
$ws.Range("H121").Value = 166.66667
$ws.Range("I121").Value = 166.66667
$ws.Range("K121").Value = 500.00001
$ws.Range("M121").Value = 809.99999
$ws.Range("H122").Value = 1115.875
$ws.Range("I122").Value = 532.25
$ws.Range("K122").Value = 4790.25
$ws.Range("M122").Value = -2340.25
$ws.Range("H128").Value = 89666.336
$ws.Range("I128").Value = 89666.336
$ws.Range("K128").Value = 268999.008
$ws.Range("M128").Value = -264019.008
$ws.Range("H129").Value = 1390.1578
$ws.Range("I129").Value = 1428
$ws.Range("J129").Value = 1348.1111
$ws.Range("K129").Value = 4284
$ws.Range("L129").Value = 4044.3333
$ws.Range("M129").Value = 716
$ws.Range("N129").Value = -14044.3333
$ws.Range("H131").Value = 2626.13
$ws.Range("I131").Value = 250
$ws.Range("J131").Value = 2657.3948
$ws.Range("K131").Value = 750
$ws.Range("L131").Value = 7972.1844
$ws.Range("M131").Value = 4290
$ws.Range("N131").Value = -18052.1844
$ws.Range("H132").Value = 997.1667
$ws.Range("I132").Value = 994.3333
$ws.Range("K132").Value = 8948.9997
$ws.Range("M132").Value = -6418.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 17836.666
$ws.Range("I41").Value = 25255
$ws.Range("J41").Value = 3000
$ws.Range("K41").Value = 25255
$ws.Range("L41").Value = 3000
$ws.Range("M41").Value = -24900
$ws.Range("N41").Value = -3710
$ws.Range("H102").Value = 1814
$ws.Range("I102").Value = 1309.3334
$ws.Range("K102").Value = 1309.3334
$ws.Range("M102").Value = 312.6666
$ws.Range("H132").Value = 3545.4849
$ws.Range("I132").Value = 3288.0386
$ws.Range("J132").Value = 4501.7144
$ws.Range("K132").Value = 9864.1158
$ws.Range("L132").Value = 13505.1432
$ws.Range("M132").Value = -7334.1158
$ws.Range("N132").Value = -18565.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 455.77777
$ws.Range("I16").Value = 455.77777
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 455.77777
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -285.77777
$ws.Range("N16").ClearContents()
$ws.Range("H100").Value = 2499.4194
$ws.Range("I100").Value = 2134.2354
$ws.Range("J100").Value = 2942.8572
$ws.Range("K100").Value = 2134.2354
$ws.Range("L100").Value = 2942.8572
$ws.Range("M100").Value = -1593.2354
$ws.Range("N100").Value = -4024.8572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 33714.2
$ws.Range("J123").Value = 33714.2
$ws.Range("L123").Value = 33714.2
$ws.Range("N123").Value = -43514.2
$ws.Range("H126").Value = 78369.30499999999
$ws.Range("I126").Value = 167983.5
$ws.Range("J126").Value = 1557.1428
$ws.Range("K126").Value = 503950.5
$ws.Range("L126").Value = 4671.428400000001
$ws.Range("M126").Value = -501480.5
$ws.Range("N126").Value = -9611.428400000001
